# Regenerate merged AHB files
# - Rename the "_old" / "_new" header suffixes to "_FV2210" / "_FV2304"
# - Turn the used range into an Excel Table (ListObject)
# - Freeze the header row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header row (row 1) cells ------------------------------------
$headers = @{
    "A1" = "Segmentname_FV2210"
    "B1" = "Segmentgruppe_FV2210"
    "C1" = "Segment_FV2210"
    "D1" = "Datenelement_FV2210"
    "E1" = "Segment ID_FV2210"
    "F1" = "Code_FV2210"
    "G1" = "Qualifier_FV2210"
    "H1" = "Beschreibung_FV2210"
    "I1" = "Bedingungsausdruck_FV2210"
    "J1" = "Bedingung_FV2210"
    "K1" = "diff"
    "L1" = "Segmentname_FV2304"
    "M1" = "Segmentgruppe_FV2304"
    "N1" = "Segment_FV2304"
    "O1" = "Datenelement_FV2304"
    "P1" = "Segment ID_FV2304"
    "Q1" = "Code_FV2304"
    "R1" = "Qualifier_FV2304"
    "S1" = "Beschreibung_FV2304"
    "T1" = "Bedingungsausdruck_FV2304"
    "U1" = "Bedingung_FV2304"
}

foreach ($addr in $headers.Keys) {
    $ws.Range($addr).Value = $headers[$addr]
}

# --- 2. Convert the used range A1:U80 into a table --------------------------
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U80"), 0, 1)
$tbl.Name = "Table1"

# --- 3. Freeze the header row ------------------------------------------------
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
